$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '44.118.47'
Set-TextValue $ws.Range('E2') '  -4.31%  '

Set-TextValue $ws.Range('D3') '2.668.88'
Set-TextValue $ws.Range('E3') '  +3.26%  '

Set-TextValue $ws.Range('D4') '0.999'
Set-TextValue $ws.Range('E4') '  -0.21%  '

Set-TextValue $ws.Range('D5') '305.52'
Set-TextValue $ws.Range('E5') '  +0.28%  '

Set-TextValue $ws.Range('D6') '97.04'
Set-TextValue $ws.Range('E6') '  -2.07%  '

Set-TextValue $ws.Range('D7') '0.588'
Set-TextValue $ws.Range('E7') '  -1.45%  '

Set-TextValue $ws.Range('E8') '  -0.03%  '

Set-TextValue $ws.Range('D9') '0.568'
Set-TextValue $ws.Range('E9') '  -1.25%  '

Set-TextValue $ws.Range('D10') '37.66'
Set-TextValue $ws.Range('E10') '  -2.83%  '

Set-TextValue $ws.Range('D11') '0.0823'
Set-TextValue $ws.Range('E11') '  -1.57%  '

Set-TextValue $ws.Range('D12') '7.93'
Set-TextValue $ws.Range('E12') '  -2.64%  '

Set-TextValue $ws.Range('D13') '3.075.69'
Set-TextValue $ws.Range('E13') '  +2.77%  '

Set-TextValue $ws.Range('E14') '  +1.54%  '

Set-TextValue $ws.Range('D15') '2.646.94'
Set-TextValue $ws.Range('E15') '  +0.35%  '

Set-TextValue $ws.Range('D16') '0.909'
Set-TextValue $ws.Range('E16') '  -0.14%  '

Set-TextValue $ws.Range('D17') '14.86'
Set-TextValue $ws.Range('E17') '  +0.34%  '

Set-TextValue $ws.Range('D18') '44.196.26'
Set-TextValue $ws.Range('E18') '  -4.60%  '

Set-TextValue $ws.Range('D19') '6.84'
Set-TextValue $ws.Range('E19') '  +3.54%  '

Set-TextValue $ws.Range('D20') '0.0₃0993'
Set-TextValue $ws.Range('E20') '  -0.91%  '

Set-TextValue $ws.Range('D21') '12.54'
Set-TextValue $ws.Range('E21') '  -2.42%  '

Set-TextValue $ws.Range('D22') '74.63'
Set-TextValue $ws.Range('E22') '  +4.51%  '

Set-TextValue $ws.Range('D23') '274.96'
Set-TextValue $ws.Range('E23') '  +1.33%  '

Set-TextValue $ws.Range('D24') '2.28'
Set-TextValue $ws.Range('E24') '  +5.95%  '

Set-TextValue $ws.Range('E25') '  +0.19%  '

Set-TextValue $ws.Range('D26') '30.32'
Set-TextValue $ws.Range('E26') '  +1.71%  '

Set-TextValue $ws.Range('D27') '1.00'
Set-TextValue $ws.Range('E27') '  +0.09%  '

Set-TextValue $ws.Range('D28') '10.39'
Set-TextValue $ws.Range('E28') '  -0.74%  '

Set-TextValue $ws.Range('E29') '  -1.98%  '

Set-TextValue $ws.Range('D30') '37.73'
Set-TextValue $ws.Range('E30') '  -3.31%  '

Set-TextValue $ws.Range('D31') '6.11'
Set-TextValue $ws.Range('E31') '  -0.05%  '

Set-TextValue $ws.Range('D32') '3.72'
Set-TextValue $ws.Range('E32') '  +3.45%  '

Set-TextValue $ws.Range('E33') '  +7.34%  '

Set-TextValue $ws.Range('D34') '153.53'
Set-TextValue $ws.Range('E34') '  +2.62%  '

Set-TextValue $ws.Range('D35') '2.80'
Set-TextValue $ws.Range('E35') '  -2.23%  '

Set-TextValue $ws.Range('D36') '0.0826'
Set-TextValue $ws.Range('E36') '  -0.95%  '

Set-TextValue $ws.Range('E37') '  -1.94%  '

Set-TextValue $ws.Range('D38') '25.19'
Set-TextValue $ws.Range('E38') '  +10.36%  '

Set-TextValue $ws.Range('D39') '0.122'
Set-TextValue $ws.Range('E39') '  +0.41%  '

Set-TextValue $ws.Range('D40') '15.85'
Set-TextValue $ws.Range('E40') '  +0.50%  '

Set-TextValue $ws.Range('D41') '3.56'
Set-TextValue $ws.Range('E41') '  -0.18%  '

Set-TextValue $ws.Range('D42') '0.0320'
Set-TextValue $ws.Range('E42') '  -2.15%  '

Set-TextValue $ws.Range('E43') '  -3.91%  '

Set-TextValue $ws.Range('D44') '2.126.61'
Set-TextValue $ws.Range('E44') '  -0.67%  '

Set-TextValue $ws.Range('D45') '0.998'
Set-TextValue $ws.Range('E45') '  -0.02%  '

Set-TextValue $ws.Range('D46') '91.43'
Set-TextValue $ws.Range('E46') '  -2.23%  '

Set-TextValue $ws.Range('D47') '9.26'
Set-TextValue $ws.Range('E47') '  -3.14%  '

Set-TextValue $ws.Range('B48') 'RocketPoolETH'
Set-TextValue $ws.Range('C48') 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
Set-TextValue $ws.Range('D48') '2.923.48'
Set-TextValue $ws.Range('E48') '  +2.97%  '

Set-TextValue $ws.Range('B49') 'Aave'
Set-TextValue $ws.Range('C49') 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws.Range('D49') '109.60'
Set-TextValue $ws.Range('E49') '  +1.09%  '

Set-TextValue $ws.Range('D50') '1.60'
Set-TextValue $ws.Range('E50') '  +3.80%  '

Set-TextValue $ws.Range('D51') '0.195'
Set-TextValue $ws.Range('E51') '  -1.09%  '
